$wb = $excel.ActiveWorkbook

# sheet1 (Worksheets.Item(1))
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 680
$ws.Range("F3").Value = 1503
$ws.Range("F4").Value = 3259
$ws.Range("F6").Value = 680
$ws.Range("F7").Value = 2248
$ws.Range("F8").Value = 487
$ws.Range("F10").Value = 237
$ws.Range("F11").Value = 130
$ws.Range("F12").Value = 324
$ws.Range("F13").Value = 1075
$ws.Range("F14").Value = 439
$ws.Range("F15").Value = 12
$ws.Range("F17").Value = 222
$ws.Range("F18").Value = 4528
$ws.Range("F19").Value = 7
$ws.Range("F20").Value = 1311
$ws.Range("F21").Value = 3437
$ws.Range("F22").Value = 327
$ws.Range("F23").Value = 101
$ws.Range("F25").Value = 3573
$ws.Range("F26").Value = 4998
$ws.Range("F28").Value = 973
$ws.Range("F29").Value = 547
$ws.Range("F30").Value = 3221
$ws.Range("F31").Value = 354
$ws.Range("F33").Value = 133
$ws.Range("F34").Value = 88
$ws.Range("F35").Value = 877
$ws.Range("F36").Value = 1166
$ws.Range("F37").Value = 1412
$ws.Range("F38").Value = 120
$ws.Range("F39").Value = 1343
$ws.Range("F40").Value = 854
$ws.Range("F41").Value = 15
$ws.Range("F42").Value = 811
$ws.Range("F43").Value = 497
$ws.Range("F44").Value = 53
$ws.Range("F45").Value = 303
$ws.Range("F46").Value = 64
$ws.Range("F47").Value = 155
$ws.Range("F48").Value = 365
$ws.Range("F49").Value = 3716

# sheet2 (Worksheets.Item(2))
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 1002
$ws.Range("F8").Value = 36
$ws.Range("F23").Value = 15

# sheet3 (Worksheets.Item(3))
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2157

# sheet4 (Worksheets.Item(4))
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2157
$ws.Range("F3").Value = 680
$ws.Range("F4").Value = 1503
$ws.Range("F5").Value = 3259
$ws.Range("F7").Value = 680
$ws.Range("F9").Value = 2248
$ws.Range("F10").Value = 487
$ws.Range("F12").Value = 238
$ws.Range("F13").Value = 1002
$ws.Range("F14").Value = 130
$ws.Range("F15").Value = 324
$ws.Range("F16").Value = 1075
$ws.Range("F17").Value = 439
$ws.Range("F18").Value = 12
$ws.Range("F19").Value = 222
$ws.Range("F20").Value = 4528
$ws.Range("F21").Value = 1311
$ws.Range("F23").Value = 3437
$ws.Range("F24").Value = 3573
$ws.Range("F25").Value = 4998
$ws.Range("F27").Value = 973
$ws.Range("F28").Value = 3221
$ws.Range("F29").Value = 354
$ws.Range("F31").Value = 133
$ws.Range("F32").Value = 88
$ws.Range("F33").Value = 877
$ws.Range("F34").Value = 1166
$ws.Range("F35").Value = 1412
$ws.Range("F36").Value = 120
$ws.Range("F37").Value = 1343
$ws.Range("F39").Value = 854
$ws.Range("F40").Value = 497
$ws.Range("F42").Value = 53
$ws.Range("F44").Value = 303
$ws.Range("F45").Value = 15
$ws.Range("F46").Value = 64
$ws.Range("F47").Value = 155
$ws.Range("F48").Value = 365
$ws.Range("F49").Value = 3716
